# Append: 2025-09-13 01:27 JST
# Update the "取得日時" (acquisition timestamp) column A on the ランサーズ sheet
# for existing data rows (2-18) from "2025-09-13 01:10:51" to "2025-09-13 01:27:18".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newValue = "2025-09-13 01:27:18"

for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 1).Value = $newValue
}
